$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# "Eintragungsdatum" (registration date) value updated from 15.12.2023 to 17.12.2023
$ws.Range("B7").Value = "17.12.2023"

# Move/confirm the active selection, matching the saved cursor position
[void]$ws.Range("B8").Select()
